# [Kadastro App] Yeni kayit eklendi: 2989
# Appends the new record row (A55:F55) to both the "Kayitlar" summary sheet
# and the matching district sheet ("Erdemli"), mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")

$rowValues = @{
    "A" = "2989"
    "B" = "2025-09-10"
    "C" = "Erdemli"
    "D" = "1"
    "E" = "ÇAP"
    "F" = "SEVİL SARAÇER (Tekniker)"
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # The sheet stores every value (including number-/date-looking ones) as
    # plain text, so force text entry with a leading apostrophe, then clear
    # the resulting "quote prefix" look so the new cells keep the sheet's
    # default (unstyled) appearance, same as every other data row.
    $ws.Range("A55").Value = "'" + $rowValues["A"]
    $ws.Range("B55").Value = "'" + $rowValues["B"]
    $ws.Range("C55").Value = $rowValues["C"]
    $ws.Range("D55").Value = "'" + $rowValues["D"]
    $ws.Range("E55").Value = $rowValues["E"]
    $ws.Range("F55").Value = $rowValues["F"]

    $ws.Range("A55:F55").Style = "Normal"
}
